$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '57.511.93'),
    @('E2', '  -4.41%  '),
    @('D3', '2.948.10'),
    @('E3', '  -1.00%  '),
    @('E4', '  +0.14%  '),
    @('D5', '555.35'),
    @('E5', '  -3.83%  '),
    @('D6', '131.99'),
    @('E6', '  +5.15%  '),
    @('E7', '  +0.16%  '),
    @('D8', '0.514'),
    @('E8', '  +2.80%  '),
    @('D9', '2.942.96'),
    @('E9', '  -1.14%  '),
    @('D10', '0.127'),
    @('E10', '  -3.98%  '),
    @('E11', '  -5.38%  '),
    @('D12', '0.448'),
    @('E12', '  +1.64%  '),
    @('E13', '  -0.08%  '),
    @('D14', '32.74'),
    @('E14', '  +0.87%  '),
    @('E15', '  +0.87%  '),
    @('D16', '3.434.39'),
    @('D17', '6.80'),
    @('E17', '  +9.71%  '),
    @('D18', '2.943.83'),
    @('E18', '  -0.85%  '),
    @('D19', '57.522.76'),
    @('E19', '  -4.20%  '),
    @('D20', '415.68'),
    @('E20', '  -3.99%  '),
    @('D21', '13.14'),
    @('E21', '  +0.49%  '),
    @('D22', '0.682'),
    @('E22', '  +3.26%  '),
    @('D23', '6.95'),
    @('E23', '  -0.85%  '),
    @('D24', '13.00'),
    @('E24', '  +2.64%  '),
    @('D25', '79.19'),
    @('E25', '  +0.20%  '),
    @('E26', '  -0.09%  '),
    @('D27', '1.00'),
    @('E27', '  +0.20%  '),
    @('E28', '  -1.62%  '),
    @('D29', '7.53'),
    @('E29', '  +4.02%  '),
    @('D30', '1.97'),
    @('E30', '  +5.01%  '),
    @('B31', 'EthereumClassic'),
    @('C31', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @('D31', '25.12'),
    @('E31', '  -0.64%  '),
    @('B32', 'NEARProtocol'),
    @('C32', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @('D32', '6.05'),
    @('E32', '  -1.56%  '),
    @('E33', '  +8.71%  '),
    @('D34', '5.63'),
    @('E34', '  +0.69%  '),
    @('B35', 'Stacks'),
    @('C35', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'),
    @('D35', '2.11'),
    @('E35', '  -2.88%  '),
    @('B36', 'Mantle'),
    @('C36', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D36', '0.939'),
    @('E36', '  -0.79%  '),
    @('D37', '48.51'),
    @('E37', '  -1.82%  '),
    @('E38', '  +4.23%  '),
    @('D39', '8.45'),
    @('E39', '  +6.61%  '),
    @('D40', '2.55'),
    @('E40', '  +3.64%  '),
    @('D41', '0.108'),
    @('E41', '  -0.44%  '),
    @('B42', 'VeChain'),
    @('C42', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D42', '0.0349'),
    @('E42', '  -2.79%  '),
    @('B43', 'Bittensor'),
    @('C43', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'),
    @('D43', '378.63'),
    @('E43', '  -1.28%  '),
    @('D44', '2.659.69'),
    @('E44', '  +1.17%  '),
    @('D46', '0.240'),
    @('E46', '  +2.00%  '),
    @('D47', '122.81'),
    @('E47', '  +3.22%  '),
    @('E48', '  +2.22%  '),
    @('D49', '1.99'),
    @('E49', '  +0.28%  '),
    @('D50', '23.40'),
    @('E50', '  +0.18%  '),
    @('D51', '2.01'),
    @('E51', '  +0.11%  ')
)

foreach ($pair in $updates) {
    $ref = $pair[0]
    $val = $pair[1]
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}